$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    # Force the cell to store as TEXT (matching the source data, which keeps
    # numeric-looking values like "12" as literal strings) while keeping the
    # cell's effective style identical to the sheet's default ("Normal") —
    # i.e. avoid leaving a stray custom number-format style behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
    $range.HorizontalAlignment = 1
    $range.VerticalAlignment = -4107
}

Set-TextCell $ws.Range("A3") "Ths"
Set-TextCell $ws.Range("B3") "12"
Set-TextCell $ws.Range("C3") "231"
Set-TextCell $ws.Range("D3") "2312"

Set-TextCell $ws.Range("A4") "Iphone"
Set-TextCell $ws.Range("B4") "100"
Set-TextCell $ws.Range("C4") "1"
Set-TextCell $ws.Range("D4") "102000"

Set-TextCell $ws.Range("A5") "Television"
Set-TextCell $ws.Range("B5") "312"
Set-TextCell $ws.Range("C5") "2132"
Set-TextCell $ws.Range("D5") "231"
